$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new header row (columns A and B)
$ws.Range("A1").Value = "Matcher"
$ws.Range("B1").Value = "Target URL"

# Set new data rows, shifting old B/C columns into A/B
$ws.Range("A2").Value = "/old-page"
$ws.Range("B2").Value = "https://example.com/new-page"

$ws.Range("A3").Value = "/legacy-section"
$ws.Range("B3").Value = "https://example.com/modern-section"

# Remove now-unused columns C:F entirely (clear in place, no shifting)
$ws.Range("C1:F3").ClearContents()
